$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K23: "Oui" -> "Non" ---
$ws.Range("K23").Value = "Non"

# --- Row 26: new journal entry (use row 24 as a style/format template) ---
$ws.Range("B24:K24").Copy($ws.Range("B26:K26"))
$ws.Range("B26").Value = 44265
$ws.Range("C26").Value = 0.72916666666666663
$ws.Range("D26").Value = 0.75
$ws.Range("E26").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure fin]]),"""",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure début]])"
$ws.Range("F26").Value = "Ma-20"
$ws.Range("G26").Value = "Code"
$ws.Range("H26").Value = "Jeu"
$ws.Range("I26").Value = "Maison"
$ws.Range("J26").Value = "J'ai améliorer les fonctions qui demande les postion ou l'on veut tirer, j'ai aussi mis en page le code"
$ws.Range("K26").Value = "non"
$ws.Rows(26).RowHeight = 43.2

# --- Row 27: new journal entry (use row 25 as a style/format template) ---
$ws.Range("B25:K25").Copy($ws.Range("B27:K27"))
$ws.Range("B27").Value = 44265
$ws.Range("C27").Value = 0.75
$ws.Range("D27").Value = 0.78125
$ws.Range("E27").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure fin]]),"""",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure début]])"
$ws.Range("F27").Value = "Ma-20"
$ws.Range("G27").Value = "Code"
$ws.Range("H27").Value = "Jeu"
$ws.Range("I27").Value = "Maison"
$ws.Range("J27").Value = "J'ai fait que lorsque qu'on tire sur une case la grille l affiche"
$ws.Range("K27").Value = "oui"
$ws.Rows(27).RowHeight = 28.8

# --- sheet view: scroll position & selection ---
$null = $ws.Range("H30").Select()
